# Apply 2020-08-18 data refresh to Fonds de solidarite volet 2 (regional x NAF) sheet.
# For each updated row, nombre_aides (col C) and montant_total (col D) are revised.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("C3").Value = "82"
$ws.Range("D3").Value = "296198.00"
$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("D2").Style

$ws.Range("C5").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("C5").Value = "131"
$ws.Range("D5").Value = "361386.40"
$ws.Range("C5").Style = $ws.Range("C4").Style
$ws.Range("D5").Style = $ws.Range("D4").Style

$ws.Range("C6").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("C6").Value = "389"
$ws.Range("D6").Value = "1030010.82"
$ws.Range("C6").Style = $ws.Range("C5").Style
$ws.Range("D6").Style = $ws.Range("D5").Style

$ws.Range("C7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("C7").Value = "78"
$ws.Range("D7").Value = "187289.00"
$ws.Range("C7").Style = $ws.Range("C6").Style
$ws.Range("D7").Style = $ws.Range("D6").Style

$ws.Range("C8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("C8").Value = "773"
$ws.Range("D8").Value = "2891058.81"
$ws.Range("C8").Style = $ws.Range("C7").Style
$ws.Range("D8").Style = $ws.Range("D7").Style

$ws.Range("C9").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("C9").Value = "27"
$ws.Range("D9").Value = "80600.00"
$ws.Range("C9").Style = $ws.Range("C8").Style
$ws.Range("D9").Style = $ws.Range("D8").Style

$ws.Range("C11").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("C11").Value = "36"
$ws.Range("D11").Value = "90177.00"
$ws.Range("C11").Style = $ws.Range("C10").Style
$ws.Range("D11").Style = $ws.Range("D10").Style

$ws.Range("C12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("C12").Value = "155"
$ws.Range("D12").Value = "464316.18"
$ws.Range("C12").Style = $ws.Range("C11").Style
$ws.Range("D12").Style = $ws.Range("D11").Style

$ws.Range("C13").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("C13").Value = "80"
$ws.Range("D13").Value = "201300.00"
$ws.Range("C13").Style = $ws.Range("C12").Style
$ws.Range("D13").Style = $ws.Range("D12").Style

$ws.Range("C16").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("C16").Value = "127"
$ws.Range("D16").Value = "572717.26"
$ws.Range("C16").Style = $ws.Range("C15").Style
$ws.Range("D16").Style = $ws.Range("D15").Style

$ws.Range("C17").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("C17").Value = "175"
$ws.Range("D17").Value = "399089.87"
$ws.Range("C17").Style = $ws.Range("C16").Style
$ws.Range("D17").Style = $ws.Range("D16").Style

$ws.Range("C35").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("C35").Value = "166"
$ws.Range("D35").Value = "471408.00"
$ws.Range("C35").Style = $ws.Range("C34").Style
$ws.Range("D35").Style = $ws.Range("D34").Style

$ws.Range("C37").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("C37").Value = "373"
$ws.Range("D37").Value = "1464941.10"
$ws.Range("C37").Style = $ws.Range("C36").Style
$ws.Range("D37").Style = $ws.Range("D36").Style

$ws.Range("C43").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("C43").Value = "28"
$ws.Range("D43").Value = "68971.00"
$ws.Range("C43").Style = $ws.Range("C42").Style
$ws.Range("D43").Style = $ws.Range("D42").Style

$ws.Range("C45").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("C45").Value = "47"
$ws.Range("D45").Value = "164048.92"
$ws.Range("C45").Style = $ws.Range("C44").Style
$ws.Range("D45").Style = $ws.Range("D44").Style

$ws.Range("C75").NumberFormat = "@"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("C75").Value = "42"
$ws.Range("D75").Value = "126579.25"
$ws.Range("C75").Style = $ws.Range("C74").Style
$ws.Range("D75").Style = $ws.Range("D74").Style

$ws.Range("C78").NumberFormat = "@"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("C78").Value = "209"
$ws.Range("D78").Value = "583693.00"
$ws.Range("C78").Style = $ws.Range("C77").Style
$ws.Range("D78").Style = $ws.Range("D77").Style

$ws.Range("C80").NumberFormat = "@"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("C80").Value = "481"
$ws.Range("D80").Value = "2076939.03"
$ws.Range("C80").Style = $ws.Range("C79").Style
$ws.Range("D80").Style = $ws.Range("D79").Style

$ws.Range("C81").NumberFormat = "@"
$ws.Range("D81").NumberFormat = "@"
$ws.Range("C81").Value = "13"
$ws.Range("D81").Value = "29000.00"
$ws.Range("C81").Style = $ws.Range("C80").Style
$ws.Range("D81").Style = $ws.Range("D80").Style

$ws.Range("C84").NumberFormat = "@"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("C84").Value = "71"
$ws.Range("D84").Value = "243572.36"
$ws.Range("C84").Style = $ws.Range("C83").Style
$ws.Range("D84").Style = $ws.Range("D83").Style

$ws.Range("C86").NumberFormat = "@"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("C86").Value = "44"
$ws.Range("D86").Value = "101500.00"
$ws.Range("C86").Style = $ws.Range("C85").Style
$ws.Range("D86").Style = $ws.Range("D85").Style

$ws.Range("C89").NumberFormat = "@"
$ws.Range("D89").NumberFormat = "@"
$ws.Range("C89").Value = "107"
$ws.Range("D89").Value = "277620.00"
$ws.Range("C89").Style = $ws.Range("C88").Style
$ws.Range("D89").Style = $ws.Range("D88").Style

$ws.Range("C107").NumberFormat = "@"
$ws.Range("D107").NumberFormat = "@"
$ws.Range("C107").Value = "71"
$ws.Range("D107").Value = "176310.00"
$ws.Range("C107").Style = $ws.Range("C106").Style
$ws.Range("D107").Style = $ws.Range("D106").Style

$ws.Range("C108").NumberFormat = "@"
$ws.Range("D108").NumberFormat = "@"
$ws.Range("C108").Value = "35"
$ws.Range("D108").Value = "115434.00"
$ws.Range("C108").Style = $ws.Range("C107").Style
$ws.Range("D108").Style = $ws.Range("D107").Style

$ws.Range("C110").NumberFormat = "@"
$ws.Range("D110").NumberFormat = "@"
$ws.Range("C110").Value = "84"
$ws.Range("D110").Value = "528606.82"
$ws.Range("C110").Style = $ws.Range("C109").Style
$ws.Range("D110").Style = $ws.Range("D109").Style

$ws.Range("C113").NumberFormat = "@"
$ws.Range("D113").NumberFormat = "@"
$ws.Range("C113").Value = "25"
$ws.Range("D113").Value = "74767.00"
$ws.Range("C113").Style = $ws.Range("C112").Style
$ws.Range("D113").Style = $ws.Range("D112").Style

$ws.Range("C117").NumberFormat = "@"
$ws.Range("D117").NumberFormat = "@"
$ws.Range("C117").Value = "19"
$ws.Range("D117").Value = "99068.92"
$ws.Range("C117").Style = $ws.Range("C116").Style
$ws.Range("D117").Style = $ws.Range("D116").Style

$ws.Range("C119").NumberFormat = "@"
$ws.Range("D119").NumberFormat = "@"
$ws.Range("C119").Value = "12"
$ws.Range("D119").Value = "32000.00"
$ws.Range("C119").Style = $ws.Range("C118").Style
$ws.Range("D119").Style = $ws.Range("D118").Style

$ws.Range("C122").NumberFormat = "@"
$ws.Range("D122").NumberFormat = "@"
$ws.Range("C122").Value = "245"
$ws.Range("D122").Value = "673508.00"
$ws.Range("C122").Style = $ws.Range("C121").Style
$ws.Range("D122").Style = $ws.Range("D121").Style

$ws.Range("C123").NumberFormat = "@"
$ws.Range("D123").NumberFormat = "@"
$ws.Range("C123").Value = "110"
$ws.Range("D123").Value = "288081.45"
$ws.Range("C123").Style = $ws.Range("C122").Style
$ws.Range("D123").Style = $ws.Range("D122").Style

$ws.Range("C124").NumberFormat = "@"
$ws.Range("D124").NumberFormat = "@"
$ws.Range("C124").Value = "476"
$ws.Range("D124").Value = "2119432.06"
$ws.Range("C124").Style = $ws.Range("C123").Style
$ws.Range("D124").Style = $ws.Range("D123").Style

$ws.Range("C132").NumberFormat = "@"
$ws.Range("D132").NumberFormat = "@"
$ws.Range("C132").Value = "85"
$ws.Range("D132").Value = "381163.75"
$ws.Range("C132").Style = $ws.Range("C131").Style
$ws.Range("D132").Style = $ws.Range("D131").Style

$ws.Range("C133").NumberFormat = "@"
$ws.Range("D133").NumberFormat = "@"
$ws.Range("C133").Value = "117"
$ws.Range("D133").Value = "296136.44"
$ws.Range("C133").Style = $ws.Range("C132").Style
$ws.Range("D133").Style = $ws.Range("D132").Style

$ws.Range("C198").NumberFormat = "@"
$ws.Range("D198").NumberFormat = "@"
$ws.Range("C198").Value = "32"
$ws.Range("D198").Value = "111574.12"
$ws.Range("C198").Style = $ws.Range("C197").Style
$ws.Range("D198").Style = $ws.Range("D197").Style

$ws.Range("C199").NumberFormat = "@"
$ws.Range("D199").NumberFormat = "@"
$ws.Range("C199").Value = "653"
$ws.Range("D199").Value = "2467424.58"
$ws.Range("C199").Style = $ws.Range("C198").Style
$ws.Range("D199").Style = $ws.Range("D198").Style

$ws.Range("C203").NumberFormat = "@"
$ws.Range("D203").NumberFormat = "@"
$ws.Range("C203").Value = "155"
$ws.Range("D203").Value = "479133.00"
$ws.Range("C203").Style = $ws.Range("C202").Style
$ws.Range("D203").Style = $ws.Range("D202").Style

$ws.Range("C207").NumberFormat = "@"
$ws.Range("D207").NumberFormat = "@"
$ws.Range("C207").Value = "124"
$ws.Range("D207").Value = "582888.14"
$ws.Range("C207").Style = $ws.Range("C206").Style
$ws.Range("D207").Style = $ws.Range("D206").Style

$ws.Range("C239").NumberFormat = "@"
$ws.Range("D239").NumberFormat = "@"
$ws.Range("C239").Value = "23"
$ws.Range("D239").Value = "67250.00"
$ws.Range("C239").Style = $ws.Range("C238").Style
$ws.Range("D239").Style = $ws.Range("D238").Style

$ws.Range("C240").NumberFormat = "@"
$ws.Range("D240").NumberFormat = "@"
$ws.Range("C240").Value = "83"
$ws.Range("D240").Value = "230538.00"
$ws.Range("C240").Style = $ws.Range("C239").Style
$ws.Range("D240").Style = $ws.Range("D239").Style

$ws.Range("C241").NumberFormat = "@"
$ws.Range("D241").NumberFormat = "@"
$ws.Range("C241").Value = "153"
$ws.Range("D241").Value = "396200.00"
$ws.Range("C241").Style = $ws.Range("C240").Style
$ws.Range("D241").Style = $ws.Range("D240").Style

$ws.Range("C242").NumberFormat = "@"
$ws.Range("D242").NumberFormat = "@"
$ws.Range("C242").Value = "510"
$ws.Range("D242").Value = "1329575.83"
$ws.Range("C242").Style = $ws.Range("C241").Style
$ws.Range("D242").Style = $ws.Range("D241").Style

$ws.Range("C243").NumberFormat = "@"
$ws.Range("D243").NumberFormat = "@"
$ws.Range("C243").Value = "98"
$ws.Range("D243").Value = "287827.11"
$ws.Range("C243").Style = $ws.Range("C242").Style
$ws.Range("D243").Style = $ws.Range("D242").Style

$ws.Range("C244").NumberFormat = "@"
$ws.Range("D244").NumberFormat = "@"
$ws.Range("C244").Value = "988"
$ws.Range("D244").Value = "3592661.86"
$ws.Range("C244").Style = $ws.Range("C243").Style
$ws.Range("D244").Style = $ws.Range("D243").Style

$ws.Range("C245").NumberFormat = "@"
$ws.Range("D245").NumberFormat = "@"
$ws.Range("C245").Value = "40"
$ws.Range("D245").Value = "103500.00"
$ws.Range("C245").Style = $ws.Range("C244").Style
$ws.Range("D245").Style = $ws.Range("D244").Style

$ws.Range("C247").NumberFormat = "@"
$ws.Range("D247").NumberFormat = "@"
$ws.Range("C247").Value = "80"
$ws.Range("D247").Value = "182500.00"
$ws.Range("C247").Style = $ws.Range("C246").Style
$ws.Range("D247").Style = $ws.Range("D246").Style

$ws.Range("C248").NumberFormat = "@"
$ws.Range("D248").NumberFormat = "@"
$ws.Range("C248").Value = "183"
$ws.Range("D248").Value = "566429.19"
$ws.Range("C248").Style = $ws.Range("C247").Style
$ws.Range("D248").Style = $ws.Range("D247").Style

$ws.Range("C249").NumberFormat = "@"
$ws.Range("D249").NumberFormat = "@"
$ws.Range("C249").Value = "125"
$ws.Range("D249").Value = "401693.00"
$ws.Range("C249").Style = $ws.Range("C248").Style
$ws.Range("D249").Style = $ws.Range("D248").Style

$ws.Range("C250").NumberFormat = "@"
$ws.Range("D250").NumberFormat = "@"
$ws.Range("C250").Value = "96"
$ws.Range("D250").Value = "259972.92"
$ws.Range("C250").Style = $ws.Range("C249").Style
$ws.Range("D250").Style = $ws.Range("D249").Style

$ws.Range("C252").NumberFormat = "@"
$ws.Range("D252").NumberFormat = "@"
$ws.Range("C252").Value = "125"
$ws.Range("D252").Value = "446812.14"
$ws.Range("C252").Style = $ws.Range("C251").Style
$ws.Range("D252").Style = $ws.Range("D251").Style

$ws.Range("C253").NumberFormat = "@"
$ws.Range("D253").NumberFormat = "@"
$ws.Range("C253").Value = "214"
$ws.Range("D253").Value = "478163.00"
$ws.Range("C253").Style = $ws.Range("C252").Style
$ws.Range("D253").Style = $ws.Range("D252").Style
